$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "65.971.82"
$ws.Cells.Item(2, 5).Value = "  +1.00%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.314.52"
$ws.Cells.Item(3, 5).Value = "  +0.55%  "
$ws.Cells.Item(4, 5).Value = "  +0.19%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "188.99"
$ws.Cells.Item(5, 5).Value = "  +5.52%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "556.38"
$ws.Cells.Item(6, 5).Value = "  +0.56%  "
$ws.Cells.Item(7, 5).Value = "  -0.13%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.582"
$ws.Cells.Item(8, 5).Value = "  -0.67%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "3.310.31"
$ws.Cells.Item(9, 5).Value = "  +0.78%  "
$ws.Cells.Item(10, 5).Value = "  -1.92%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.583"
$ws.Cells.Item(11, 5).Value = "  +0.28%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "46.99"
$ws.Cells.Item(12, 5).Value = "  +0.10%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000268"
$ws.Cells.Item(13, 5).Value = "  +2.43%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "8.65"
$ws.Cells.Item(14, 5).Value = "  +1.88%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.844.14"
$ws.Cells.Item(15, 5).Value = "  +0.72%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "601.85"
$ws.Cells.Item(16, 5).Value = "  -0.21%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "65.978.99"
$ws.Cells.Item(17, 5).Value = "  +1.13%  "
$ws.Cells.Item(18, 2).Value = "TRON"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.118"
$ws.Cells.Item(18, 5).Value = "  +1.01%  "
$ws.Cells.Item(19, 2).Value = "Chainlink"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "17.91"
$ws.Cells.Item(19, 5).Value = "  -0.40%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "3.312.70"
$ws.Cells.Item(20, 5).Value = "  +0.80%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "11.03"
$ws.Cells.Item(21, 5).Value = "  -2.72%  "
$ws.Cells.Item(22, 5).Value = "  +0.21%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "18.58"
$ws.Cells.Item(23, 5).Value = "  +6.61%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "5.11"
$ws.Cells.Item(24, 5).Value = "  +3.07%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "101.11"
$ws.Cells.Item(25, 5).Value = "  -1.29%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.95"
$ws.Cells.Item(26, 5).Value = "  +0.03%  "
$ws.Cells.Item(27, 2).Value = "LEO"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "6.02"
$ws.Cells.Item(27, 5).Value = "  +0.94%  "
$ws.Cells.Item(28, 2).Value = "ImmutableX"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.74"
$ws.Cells.Item(28, 5).Value = "  +3.23%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "9.52"
$ws.Cells.Item(29, 5).Value = "  +2.42%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "8.69"
$ws.Cells.Item(30, 5).Value = "  +0.77%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "30.27"
$ws.Cells.Item(31, 5).Value = "  -0.12%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "6.74"
$ws.Cells.Item(32, 5).Value = "  +8.65%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.84"
$ws.Cells.Item(33, 5).Value = "  -1.52%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "575.23"
$ws.Cells.Item(34, 5).Value = "  +8.57%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "11.02"
$ws.Cells.Item(35, 5).Value = "  +0.72%  "
$ws.Cells.Item(36, 5).Value = "  +0.53%  "
$ws.Cells.Item(37, 5).Value = "  +0.15%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "3.707.81"
$ws.Cells.Item(38, 5).Value = "  -2.34%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "57.03"
$ws.Cells.Item(39, 5).Value = "  +1.86%  "
$ws.Cells.Item(40, 2).Value = "CoreDAO"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.49"
$ws.Cells.Item(40, 5).Value = "  +9.62%  "
$ws.Cells.Item(41, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "33.78"
$ws.Cells.Item(41, 5).Value = "  +6.89%  "
$ws.Cells.Item(42, 2).Value = "Stacks"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.26"
$ws.Cells.Item(42, 5).Value = "  -4.59%  "
$ws.Cells.Item(43, 2).Value = "Kaspa"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.129"
$ws.Cells.Item(43, 5).Value = "  +2.93%  "
$ws.Cells.Item(44, 2).Value = "Fetch.AI"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.68"
$ws.Cells.Item(44, 5).Value = "  +1.79%  "
$ws.Cells.Item(45, 2).Value = "PEPE"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.0₃0701"
$ws.Cells.Item(45, 5).Value = "  -0.91%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "3.41"
$ws.Cells.Item(46, 5).Value = "  +5.49%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.339"
$ws.Cells.Item(47, 5).Value = "  +0.78%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0418"
$ws.Cells.Item(48, 5).Value = "  +2.72%  "
$ws.Cells.Item(49, 5).Value = "  +0.18%  "
$ws.Cells.Item(50, 5).Value = "  -0.22%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.00"
$ws.Cells.Item(51, 5).Value = "  +0.13%  "
